# Refresh the cryptocurrency price/volume table (and resolve the Aave/ImmutableX row-order swap)
# with the latest scraped values, mirroring the GitHub Actions data-update commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    $cell = $ws.Range($rangeAddress)
    # Prefix with an apostrophe so Excel stores numeric-looking strings
    # (e.g. "224.21") as text rather than silently coercing them to numbers,
    # then strip the resulting quote-prefix formatting so the cell keeps the
    # workbook default style (matching the original unstyled inline strings).
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

Set-TextValue "D2" "31.615.43"
Set-TextValue "E2" "  +5.98%  "
Set-TextValue "D3" "1.716.65"
Set-TextValue "E3" "  +4.81%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.17%  "
Set-TextValue "D5" "224.21"
Set-TextValue "E5" "  +4.14%  "
Set-TextValue "D6" "0.538"
Set-TextValue "E6" "  +3.68%  "
Set-TextValue "D7" "0.998"
Set-TextValue "E7" "  -0.17%  "
Set-TextValue "D8" "30.03"
Set-TextValue "E8" "  +3.96%  "
Set-TextValue "D9" "0.270"
Set-TextValue "E9" "  +3.47%  "
Set-TextValue "E10" "  +6.91%  "
Set-TextValue "D11" "0.0911"
Set-TextValue "E11" "  +1.32%  "
Set-TextValue "D12" "1.955.65"
Set-TextValue "E12" "  +4.39%  "
Set-TextValue "D13" "1.716.28"
Set-TextValue "E13" "  +4.54%  "
Set-TextValue "D14" "0.616"
Set-TextValue "E14" "  +4.38%  "
Set-TextValue "D15" "10.19"
Set-TextValue "E15" "  +7.88%  "
Set-TextValue "E16" "  +7.92%  "
Set-TextValue "D17" "31.594.50"
Set-TextValue "E17" "  +5.84%  "
Set-TextValue "D18" "67.52"
Set-TextValue "E18" "  +4.85%  "
Set-TextValue "D19" "251.74"
Set-TextValue "E19" "  +5.44%  "
Set-TextValue "E20" "  +3.17%  "
Set-TextValue "E21" "  -0.09%  "
Set-TextValue "D22" "10.15"
Set-TextValue "E22" "  +2.26%  "
Set-TextValue "D23" "4.26"
Set-TextValue "E23" "  +3.17%  "
Set-TextValue "D24" "2.18"
Set-TextValue "E24" "  -0.11%  "
Set-TextValue "D25" "159.51"
Set-TextValue "E25" "  +1.66%  "
Set-TextValue "D26" "16.11"
Set-TextValue "E26" "  +3.58%  "
Set-TextValue "E27" "  +3.56%  "
Set-TextValue "D28" "6.82"
Set-TextValue "E28" "  +3.08%  "
Set-TextValue "E29" "  -0.13%  "
Set-TextValue "D30" "3.89"
Set-TextValue "E30" "  +15.14%  "
Set-TextValue "E31" "  +1.82%  "
Set-TextValue "E32" "  +4.59%  "
Set-TextValue "D33" "3.42"
Set-TextValue "E33" "  +7.02%  "
Set-TextValue "D34" "1.533.04"
Set-TextValue "E34" "  +7.96%  "
Set-TextValue "E35" "  +3.99%  "
Set-TextValue "E36" "  +2.54%  "
Set-TextValue "B37" "Aave"
Set-TextValue "C37" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D37" "83.11"
Set-TextValue "E37" "  +8.75%  "
Set-TextValue "B38" "ImmutableX"
Set-TextValue "C38" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D38" "0.614"
Set-TextValue "E38" "  +8.74%  "
Set-TextValue "E39" "  +5.04%  "
Set-TextValue "E40" "  +0.25%  "
Set-TextValue "E41" "  +0.64%  "
Set-TextValue "D42" "2.04"
Set-TextValue "E42" "  +5.18%  "
Set-TextValue "D43" "0.855"
Set-TextValue "E43" "  +2.74%  "
Set-TextValue "E44" "  +0.88%  "
Set-TextValue "E45" "  +3.48%  "
Set-TextValue "D46" "0.999"
Set-TextValue "E46" "  -0.06%  "
Set-TextValue "D47" "52.59"
Set-TextValue "E47" "  +5.99%  "
Set-TextValue "E48" "  +5.14%  "
Set-TextValue "D49" "1.848.72"
Set-TextValue "E49" "  +3.80%  "
Set-TextValue "D50" "0.0₆0118"
Set-TextValue "E50" "  +9.09%  "
Set-TextValue "D51" "93.62"
Set-TextValue "E51" "  +0.47%  "
